# Add a new "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: same text style as the other header cells (bold/centered/bordered),
# so copy the formatting from the neighboring header cell (G1) onto H1.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data cell: numeric value for row 2.
$ws.Range("H2").Value = 0
